$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 currently holds phone "09876543" as text (leading zero kept) with
# 0 points. Insert a new row 34 that preserves that original text record
# (Copy keeps the literal text - a plain .Value assignment would coerce the
# numeric-looking string back into a number), then turn row 33's phone into
# the numeric form 9876543, with its points reset to 0.00.

$ws.Range("A33:C33").Copy($ws.Range("A34:C34"))

$ws.Range("A33").Value = 9876543
$ws.Range("C33").Value = 0
